$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 44.0
$ws.Range("C2").Value = 96.0
$ws.Range("D2").Value = 2.0
$ws.Range("E2").Value = 20.0
$ws.Range("F2").Value = 0.006619216621862723
$ws.Range("G2").Value = 0.09067422947244705
$ws.Range("H2").Value = 0.09705053943736999
$ws.Range("I2").Value = 0.03317966590840471
